$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update MACRO_SCORE (N2)
$ws.Range("N2").Value = 85.92500513438651

# Row 3: update MACRO_SCORE (N3)
$ws.Range("N3").Value = 85.92500513438651

# Row 4: update 종가(D4), RSI(E4), 5일수익률(F4), 5일상승확률(I4), 최종점수(K4), MACRO_SCORE(N4)
$ws.Range("D4").Value = 4281.8
$ws.Range("E4").Value = 73.5
$ws.Range("F4").Value = 5.43
$ws.Range("I4").Value = 56
$ws.Range("K4").Value = 60.2
$ws.Range("N4").Value = 85.92500513438651
